$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.098.47'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '2.240.38'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.573'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0813'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.581.84'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.239.36'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.828'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.57'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.35%  '
$ws.Range('D18').Value = '43.993.09'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').Value = '0.0₃0966'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '39.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0802'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.120'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.109'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0298'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '1.721.33'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '83.82'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.189'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('E48').Value = '  -4.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '69.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.33'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.97%  '
